# Update workbook "Avverkningsanmälningar" (Övertorneå logging notifications list).
#
# Summary of change:
#  1. A new record for "A 36117-2023" (with updated figures) is inserted as row 3,
#     pushing the former rows 3 and 4 (A 21840-2023, A 19650-2023) down to rows 4 and 5.
#  2. The old "A 36117-2023" row (which ends up at row 6 after the insert, with its
#     previous/smaller figures) is removed, so the old rows 6..372 slide back up and
#     keep their original row numbers (6..372).
#  3. The "Förändrad" date in column C is bumped from 2023-09-19 (45188) to
#     2023-09-20 (45189) for every data row (2..372).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at position 3 ------------------------------------
$ws.Rows("3:3").Insert()

# --- Step 2: remove the old "A 36117-2023" row, now shifted down to row 6 ------
$ws.Rows("6:6").Delete()

# --- Step 3: populate the new row 3 with the refreshed "A 36117-2023" data -----
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
$ws.Range("R3").WrapText = $true

$ws.Range("A3").Value = "A 36117-2023"
$ws.Range("B3").Value = 45148
$ws.Range("C3").Value = 45189
$ws.Range("D3").Value = "NORRBOTTENS LÄN"
$ws.Range("E3").Value = "ÖVERTORNEÅ"
$ws.Range("G3").Value = 15.5
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 7
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 14
$ws.Range("R3").Value = "Gammelgransskål`r`nGarnlav`r`nMörk kolflarnlav`r`nSpillkråka`r`nTallticka`r`nUllticka`r`nVitplätt`r`nBronshjon`r`nBårdlav`r`nSkinnlav`r`nStor aspticka`r`nStuplav`r`nVedticka`r`nRevlummer"

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/artfynd/A 36117-2023.xlsx", "A 36117-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/kartor/A 36117-2023.png", "A 36117-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/klagomål/A 36117-2023.docx", "A 36117-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/klagomålsmail/A 36117-2023.docx", "A 36117-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/tillsyn/A 36117-2023.docx", "A 36117-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVERTORNEA/tillsynsmail/A 36117-2023.docx", "A 36117-2023")'

# keep the row height consistent with the rest of the sheet (15pt), since inserting
# a heavily-wrapped cell otherwise triggers Excel's auto row-height calculation.
$ws.Rows("3:3").RowHeight = 15

# --- Step 4: bump the "Förändrad" date (column C) on every data row ------------
for ($r = 2; $r -le 372; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}
